# Rename the "Property1" worksheet to "DataNode" so the sheet naming
# reflects the unified DataNode/DataTable/Entity concept described in the
# commit message ("unify the conception of DataNode, DataTable, Entity").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# The sheet has a frozen pane at row 8 (state="frozen", ySplit=8). Move the
# saved selection/active cell of the lower-left (scrollable) pane from K9
# to O40, matching the view state recorded in the edited workbook.
$ws.Range("O40").Select()
